$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 8063.7334
$ws.Range("I40").Value = 5499.6665
$ws.Range("K40").Value = 5499.6665
$ws.Range("M40").Value = -5324.6665
# Row 43
$ws.Range("H43").Value = 4999.3335
$ws.Range("J43").Value = 3999.6667
$ws.Range("L43").Value = 3999.6667
$ws.Range("N43").Value = -4137.6667
# Row 80
$ws.Range("H80").Value = 5322.727
$ws.Range("J80").Value = 9175
$ws.Range("L80").Value = 27525
$ws.Range("N80").Value = -29521
# Row 83
$ws.Range("H83").Value = 5322.727
$ws.Range("J83").Value = 9175
$ws.Range("L83").Value = 82575
$ws.Range("N83").Value = -92559
# Row 111
$ws.Range("H111").Value = 5249
$ws.Range("J111").Value = 5998
$ws.Range("L111").Value = 17994
$ws.Range("N111").Value = -24128
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("N136").Value = 0

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 920.4
$ws.Range("J74").Value = 882.5
$ws.Range("L74").Value = 882.5
$ws.Range("N74").Value = -2630.5
# Row 77
$ws.Range("H77").Value = 920.4
$ws.Range("J77").Value = 882.5
$ws.Range("L77").Value = 4412.5
$ws.Range("N77").Value = -13148.5
# Row 121
$ws.Range("H121").Value = 99995
$ws.Range("J121").Value = 99995
$ws.Range("L121").Value = 99995
$ws.Range("N121").Value = -103489

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 99995
$ws.Range("J2").Value = 99995
$ws.Range("L2").Value = 99995
$ws.Range("N2").Value = -100221
# Row 13
$ws.Range("H13").Value = 65000
$ws.Range("J13").Value = 65000
$ws.Range("L13").Value = 65000
$ws.Range("N13").Value = -65336
# Row 22
$ws.Range("H22").Value = 85
$ws.Range("I22").Value = 85
$ws.Range("K22").Value = 85
$ws.Range("M22").Value = 88
# Row 105
$ws.Range("H105").Value = 3416.5
$ws.Range("I105").Value = 3199.8
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 3199.8
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -1452.8
$ws.Range("N105").Value = -7994
# Row 107
$ws.Range("H107").Value = 2778.7144
$ws.Range("I107").Value = 1590.2
$ws.Range("K107").Value = 1590.2
$ws.Range("M107").Value = 329.8
# Row 109
$ws.Range("H109").Value = 99995
$ws.Range("J109").Value = 99995
$ws.Range("L109").Value = 99995
$ws.Range("N109").Value = -102769
# Row 135
$ws.Range("H135").Value = 99995
$ws.Range("J135").Value = 99995
$ws.Range("L135").Value = 99995
$ws.Range("N135").Value = -110135
# Row 140
$ws.Range("H140").Value = 95780
$ws.Range("J140").Value = 95780
$ws.Range("L140").Value = 95780
$ws.Range("N140").Value = -106140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1500
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -797
# Row 74
$ws.Range("H74").Value = 46328.5
$ws.Range("J74").Value = 46328.5
$ws.Range("L74").Value = 46328.5
$ws.Range("N74").Value = -48076.5
# Row 77
$ws.Range("H77").Value = 46328.5
$ws.Range("J77").Value = 46328.5
$ws.Range("L77").Value = 138985.5
$ws.Range("N77").Value = -147721.5
# Row 122
$ws.Range("H122").Value = 1582.75
$ws.Range("J122").Value = 1663
$ws.Range("L122").Value = 4989
$ws.Range("N122").Value = -9889
# Row 127
$ws.Range("H127").Value = 99995
$ws.Range("J127").Value = 99995
$ws.Range("L127").Value = 99995
$ws.Range("N127").Value = -109915
# Row 136
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0
# Row 43
$ws.Range("H43").Value = 4468.5557
$ws.Range("I43").Value = 3152.125
$ws.Range("K43").Value = 3152.125
$ws.Range("M43").Value = -3001.125
# Row 46
$ws.Range("H46").Value = 10217
$ws.Range("I46").Value = 10217
$ws.Range("K46").Value = 10217
$ws.Range("M46").Value = -10061
# Row 119
$ws.Range("H119").Value = 99995
$ws.Range("J119").Value = 99995
$ws.Range("L119").Value = 99995
$ws.Range("N119").Value = -109671
# Row 122
$ws.Range("H122").Value = 1225
$ws.Range("I122").Value = 1225
$ws.Range("K122").Value = 3675
$ws.Range("M122").Value = -1225

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 22
$ws.Range("H22").Value = 2060.5334
$ws.Range("I22").Value = 1019
$ws.Range("J22").Value = 2581.3
$ws.Range("K22").Value = 1019
$ws.Range("L22").Value = 2581.3
$ws.Range("M22").Value = -724
$ws.Range("N22").Value = -3171.3
# Row 27
$ws.Range("H27").Value = 2060.5334
$ws.Range("I27").Value = 1019
$ws.Range("J27").Value = 2581.3
$ws.Range("K27").Value = 1019
$ws.Range("L27").Value = 2581.3
$ws.Range("M27").Value = -912
$ws.Range("N27").Value = -2795.3
# Row 46
$ws.Range("H46").Value = 4420.8667
$ws.Range("J46").Value = 4677.923
$ws.Range("L46").Value = 4677.923
$ws.Range("N46").Value = -5053.923
# Row 61
$ws.Range("H61").Value = 3500.1667
$ws.Range("I61").Value = 3500.1667
$ws.Range("K61").Value = 3500.1667
$ws.Range("M61").Value = -3298.1667
# Row 113
$ws.Range("H113").Value = 3500.1667
$ws.Range("I113").Value = 3500.1667
$ws.Range("K113").Value = 3500.1667
$ws.Range("M113").Value = -1330.1667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 30
$ws.Range("H30").Value = 10336
$ws.Range("J30").Value = 11499.5
$ws.Range("L30").Value = 11499.5
$ws.Range("N30").Value = -11713.5
# Row 80
$ws.Range("H80").Value = 6000
$ws.Range("J80").Value = 6000
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7996
# Row 83
$ws.Range("H83").Value = 6000
$ws.Range("J83").Value = 6000
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27984
# Row 100
$ws.Range("H100").Value = 376.66666
$ws.Range("I100").Value = 255.85715
$ws.Range("K100").Value = 511.7143
$ws.Range("M100").Value = 29.28570000000002

Write-Output "Applied 160 cell updates across 8 sheets."